$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F3, F4, F5
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1223
$ws1.Range("F4").Value = 2708
$ws1.Range("F5").Value = 241

# Sheet "全部类型" (sheet4): update F5, F6, F8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1223
$ws4.Range("F6").Value = 2708
$ws4.Range("F8").Value = 241
